$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-26 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-27 Saturday", 2)

$d.Content.Find.Execute("236÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "757÷5=", 2)
$d.Content.Find.Execute("620÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "771÷8=", 2)
$d.Content.Find.Execute("658÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "496÷3=", 2)
$d.Content.Find.Execute("522÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "885÷5=", 2)
$d.Content.Find.Execute("205÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "650÷3=", 2)

$d.Content.Find.Execute("794÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "732÷2=", 2)
$d.Content.Find.Execute("772÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "223÷3=", 2)
$d.Content.Find.Execute("968÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "428÷8=", 2)
$d.Content.Find.Execute("197÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "621÷5=", 2)
$d.Content.Find.Execute("645÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "928÷7=", 2)

$d.Content.Find.Execute("762÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "479÷9=", 2)
$d.Content.Find.Execute("825÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "616÷4=", 2)
$d.Content.Find.Execute("666÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "147÷5=", 2)
$d.Content.Find.Execute("693÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "403÷7=", 2)
$d.Content.Find.Execute("215÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "660÷8=", 2)

$d.Content.Find.Execute("445÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "238÷7=", 2)
$d.Content.Find.Execute("672÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "912÷5=", 2)
$d.Content.Find.Execute("700÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "787÷5=", 2)
$d.Content.Find.Execute("444÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "174÷8=", 2)
$d.Content.Find.Execute("424÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "236÷5=", 2)

$d.Content.Find.Execute("697÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "149÷3=", 2)
$d.Content.Find.Execute("480÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "462÷8=", 2)
$d.Content.Find.Execute("906÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "283÷3=", 2)
$d.Content.Find.Execute("406÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "919÷9=", 2)
$d.Content.Find.Execute("142÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "254÷7=", 2)
